$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Model'
$ws.Cells.Item(1, 2).Value = 'Mania'
$ws.Cells.Item(1, 3).Value = 'Depression'
$ws.Cells.Item(1, 4).Value = 'Schizophrenia'
$ws.Cells.Item(2, 1).Value = 'gpt-5.1 [T0-With-K]'
$ws.Cells.Item(2, 2).Value = '0.65 / 0.73'
$ws.Cells.Item(2, 3).Value = '0.76 / 0.82'
$ws.Cells.Item(2, 4).Value = '0.59 / 0.65'
$ws.Cells.Item(3, 1).Value = 'gpt-5.1 [T0.5-With-K]'
$ws.Cells.Item(3, 2).Value = '0.66 (0.65 - 0.67) / 0.76 (0.75 - 0.77)'
$ws.Cells.Item(3, 3).Value = '0.72 (0.70 - 0.76) / 0.78 (0.75 - 0.83)'
$ws.Cells.Item(3, 4).Value = '0.53 (0.49 - 0.59) / 0.60 (0.55 - 0.66)'
$ws.Cells.Item(4, 1).Value = 'gpt-5.1 [T0.5-No-K]'
$ws.Cells.Item(4, 2).Value = '0.63 (0.60 - 0.65) / 0.73 (0.70 - 0.77)'
$ws.Cells.Item(4, 3).Value = '0.70 (0.52 - 0.83) / 0.76 (0.59 - 0.87)'
$ws.Cells.Item(4, 4).Value = '0.48 (0.47 - 0.48) / 0.57 (0.56 - 0.59)'
$ws.Cells.Item(5, 1).Value = 'moonshotai_Kimi-K2-Thinking [T0.5-With-K]'
$ws.Cells.Item(5, 2).Value = '0.58 (0.52 - 0.64) / 0.68 (0.61 - 0.73)'
$ws.Cells.Item(5, 3).Value = '0.63 (0.57 - 0.70) / 0.74 (0.67 - 0.80)'
$ws.Cells.Item(5, 4).Value = '0.48 (0.46 - 0.49) / 0.58 (0.56 - 0.61)'
$ws.Cells.Item(6, 1).Value = 'gemini-2.5-flash [T0.5-No-K]'
$ws.Cells.Item(6, 2).Value = '0.64 (0.59 - 0.67) / 0.71 (0.69 - 0.72)'
$ws.Cells.Item(6, 3).Value = '0.61 (0.35 - 0.77) / 0.72 (0.45 - 0.88)'
$ws.Cells.Item(6, 4).Value = '0.48 (0.45 - 0.50) / 0.55 (0.51 - 0.57)'
$ws.Cells.Item(7, 1).Value = 'gemini-2.5-flash [T0-With-K]'
$ws.Cells.Item(7, 2).Value = '0.47 / 0.59'
$ws.Cells.Item(7, 3).Value = '0.65 / 0.74'
$ws.Cells.Item(7, 4).Value = '0.47 / 0.55'
$ws.Cells.Item(8, 1).Value = 'gemini-3-pro-preview [T0.5-No-K]'
$ws.Cells.Item(8, 2).Value = '0.60 (0.57 - 0.64) / 0.69 (0.68 - 0.71)'
$ws.Cells.Item(8, 3).Value = '0.77 (0.76 - 0.78) / 0.86 (0.83 - 0.87)'
$ws.Cells.Item(8, 4).Value = '0.47 (0.41 - 0.52) / 0.54 (0.50 - 0.58)'
$ws.Cells.Item(9, 1).Value = 'gemini-3-pro-preview [T0-No-K]'
$ws.Cells.Item(9, 2).Value = '0.60 / 0.70'
$ws.Cells.Item(9, 3).Value = '0.77 / 0.85'
$ws.Cells.Item(9, 4).Value = '0.46 / 0.55'
$ws.Cells.Item(10, 1).Value = 'gemini-3-pro-preview [T0.5-With-K]'
$ws.Cells.Item(10, 2).Value = '0.62 (0.59 - 0.64) / 0.72 (0.71 - 0.73)'
$ws.Cells.Item(10, 3).Value = '0.82 (0.81 - 0.83) / 0.87 (0.86 - 0.87)'
$ws.Cells.Item(10, 4).Value = '0.46 (0.44 - 0.47) / 0.56 (0.55 - 0.58)'
$ws.Cells.Item(11, 1).Value = 'gemini-3-pro-preview [T0-With-K]'
$ws.Cells.Item(11, 2).Value = '0.62 / 0.72'
$ws.Cells.Item(11, 3).Value = '0.83 / 0.87'
$ws.Cells.Item(11, 4).Value = '0.45 / 0.56'
$ws.Cells.Item(12, 1).Value = 'moonshotai_Kimi-K2-Thinking [T0-With-K]'
$ws.Cells.Item(12, 2).Value = '0.46 / 0.57'
$ws.Cells.Item(12, 3).Value = '0.64 / 0.74'
$ws.Cells.Item(12, 4).Value = '0.45 / 0.55'
$ws.Cells.Item(13, 1).Value = 'gpt-5.1 [T0-No-K]'
$ws.Cells.Item(13, 2).Value = '0.64 / 0.73'
$ws.Cells.Item(13, 3).Value = '0.79 / 0.86'
$ws.Cells.Item(13, 4).Value = '0.43 / 0.53'
$ws.Cells.Item(14, 1).Value = 'deepseek-ai_DeepSeek-R1 [T0-No-K]'
$ws.Cells.Item(14, 2).Value = '0.63 / 0.69'
$ws.Cells.Item(14, 3).Value = '0.79 / 0.87'
$ws.Cells.Item(14, 4).Value = '0.43 / 0.51'
$ws.Cells.Item(15, 1).Value = 'Qwen_Qwen3-Next-80B-A3B-Thinking [T0-No-K]'
$ws.Cells.Item(15, 2).Value = '0.63 / 0.69'
$ws.Cells.Item(15, 3).Value = '0.72 / 0.83'
$ws.Cells.Item(15, 4).Value = '0.42 / 0.49'
$ws.Cells.Item(16, 1).Value = 'gemini-2.5-flash [T0.5-With-K]'
$ws.Cells.Item(16, 2).Value = '0.53 (0.53 - 0.54) / 0.67 (0.66 - 0.67)'
$ws.Cells.Item(16, 3).Value = '0.64 (0.60 - 0.67) / 0.72 (0.68 - 0.75)'
$ws.Cells.Item(16, 4).Value = '0.41 (0.40 - 0.44) / 0.52 (0.50 - 0.53)'
$ws.Cells.Item(17, 1).Value = 'deepseek-ai_DeepSeek-R1 [T0.5-No-K]'
$ws.Cells.Item(17, 2).Value = '0.55 (0.48 - 0.63) / 0.61 (0.52 - 0.70)'
$ws.Cells.Item(17, 3).Value = '0.41 (0.23 - 0.75) / 0.50 (0.32 - 0.84)'
$ws.Cells.Item(17, 4).Value = '0.41 (0.31 - 0.47) / 0.50 (0.41 - 0.55)'
$ws.Cells.Item(18, 1).Value = 'claude-sonnet-4-5 [T0.5-With-K]'
$ws.Cells.Item(18, 2).Value = '0.57 (0.53 - 0.59) / 0.67 (0.65 - 0.69)'
$ws.Cells.Item(18, 3).Value = '0.61 (0.59 - 0.64) / 0.70 (0.68 - 0.73)'
$ws.Cells.Item(18, 4).Value = '0.40 (0.40 - 0.41) / 0.51 (0.51 - 0.52)'
$ws.Cells.Item(19, 1).Value = 'Qwen_Qwen3-Next-80B-A3B-Thinking [T0.5-No-K]'
$ws.Cells.Item(19, 2).Value = '0.60 (0.58 - 0.62) / 0.67 (0.65 - 0.68)'
$ws.Cells.Item(19, 3).Value = '0.73 (0.72 - 0.74) / 0.83 (0.83 - 0.84)'
$ws.Cells.Item(19, 4).Value = '0.40 (0.37 - 0.45) / 0.47 (0.44 - 0.50)'
$ws.Cells.Item(20, 1).Value = 'Qwen_Qwen3-Next-80B-A3B-Thinking [T0.5-With-K]'
$ws.Cells.Item(20, 2).Value = '0.46 (0.44 - 0.48) / 0.55 (0.52 - 0.58)'
$ws.Cells.Item(20, 3).Value = '0.46 (0.44 - 0.50) / 0.58 (0.56 - 0.61)'
$ws.Cells.Item(20, 4).Value = '0.40 (0.37 - 0.43) / 0.47 (0.45 - 0.49)'
$ws.Cells.Item(21, 1).Value = 'claude-sonnet-4-5 [T0-No-K]'
$ws.Cells.Item(21, 2).Value = '0.58 / 0.67'
$ws.Cells.Item(21, 3).Value = '0.75 / 0.85'
$ws.Cells.Item(21, 4).Value = '0.39 / 0.53'
$ws.Cells.Item(22, 1).Value = 'claude-sonnet-4-5 [T0-With-K]'
$ws.Cells.Item(22, 2).Value = '0.51 / 0.62'
$ws.Cells.Item(22, 3).Value = '0.65 / 0.75'
$ws.Cells.Item(22, 4).Value = '0.39 / 0.50'
$ws.Cells.Item(23, 1).Value = 'gemini-2.5-flash [T0-No-K]'
$ws.Cells.Item(23, 2).Value = '0.56 / 0.67'
$ws.Cells.Item(23, 3).Value = '0.72 / 0.84'
$ws.Cells.Item(23, 4).Value = '0.39 / 0.49'
$ws.Cells.Item(24, 1).Value = 'mistral-large-latest [T0-No-K]'
$ws.Cells.Item(24, 2).Value = '0.60 / 0.70'
$ws.Cells.Item(24, 3).Value = '0.73 / 0.82'
$ws.Cells.Item(24, 4).Value = '0.38 / 0.48'
$ws.Cells.Item(25, 1).Value = 'moonshotai_Kimi-K2-Thinking [T0-No-K]'
$ws.Cells.Item(25, 2).Value = '0.70 / 0.75'
$ws.Cells.Item(25, 3).Value = '0.46 / 0.56'
$ws.Cells.Item(25, 4).Value = '0.38 / 0.47'
$ws.Cells.Item(26, 1).Value = 'Qwen_Qwen3-Next-80B-A3B-Thinking [T0-With-K]'
$ws.Cells.Item(26, 2).Value = '0.47 / 0.57'
$ws.Cells.Item(26, 3).Value = '0.46 / 0.58'
$ws.Cells.Item(26, 4).Value = '0.38 / 0.46'
$ws.Cells.Item(27, 1).Value = 'openai_gpt-oss-20b [T0.5-With-K]'
$ws.Cells.Item(27, 2).Value = '0.47 (0.37 - 0.53) / 0.53 (0.43 - 0.60)'
$ws.Cells.Item(27, 3).Value = '0.43 (0.36 - 0.50) / 0.52 (0.48 - 0.59)'
$ws.Cells.Item(27, 4).Value = '0.36 (0.33 - 0.38) / 0.44 (0.41 - 0.48)'
$ws.Cells.Item(28, 1).Value = 'claude-sonnet-4-5 [T0.5-No-K]'
$ws.Cells.Item(28, 2).Value = '0.58 (0.58 - 0.58) / 0.66 (0.66 - 0.67)'
$ws.Cells.Item(28, 3).Value = '0.74 (0.74 - 0.75) / 0.84 (0.82 - 0.85)'
$ws.Cells.Item(28, 4).Value = '0.36 (0.32 - 0.39) / 0.47 (0.43 - 0.52)'
$ws.Cells.Item(29, 1).Value = 'gpt-4o-mini [T0-No-K]'
$ws.Cells.Item(29, 2).Value = '0.66 / 0.70'
$ws.Cells.Item(29, 3).Value = '0.70 / 0.82'
$ws.Cells.Item(29, 4).Value = '0.35 / 0.42'
$ws.Cells.Item(30, 1).Value = 'moonshotai_Kimi-K2-Thinking [T0.5-No-K]'
$ws.Cells.Item(30, 2).Value = '0.64 (0.60 - 0.71) / 0.71 (0.67 - 0.74)'
$ws.Cells.Item(30, 3).Value = '0.57 (0.16 - 0.78) / 0.65 (0.25 - 0.87)'
$ws.Cells.Item(30, 4).Value = '0.34 (0.27 - 0.42) / 0.44 (0.40 - 0.51)'
$ws.Cells.Item(31, 1).Value = 'openai_gpt-oss-20b [T0.5-No-K]'
$ws.Cells.Item(31, 2).Value = '0.43 (0.21 - 0.56) / 0.52 (0.29 - 0.65)'
$ws.Cells.Item(31, 3).Value = '0.15 (0.14 - 0.16) / 0.25 (0.23 - 0.27)'
$ws.Cells.Item(31, 4).Value = '0.32 (0.29 - 0.36) / 0.40 (0.34 - 0.43)'
$ws.Cells.Item(32, 1).Value = 'mistral-large-latest [T0-With-K]'
$ws.Cells.Item(32, 2).Value = '0.50 / 0.62'
$ws.Cells.Item(32, 3).Value = '0.71 / 0.77'
$ws.Cells.Item(32, 4).Value = '0.31 / 0.43'
$ws.Cells.Item(33, 1).Value = 'deepseek-ai_DeepSeek-R1 [T0-With-K]'
$ws.Cells.Item(33, 2).Value = '0.45 / 0.55'
$ws.Cells.Item(33, 3).Value = '0.35 / 0.43'
$ws.Cells.Item(33, 4).Value = '0.31 / 0.41'
$ws.Cells.Item(34, 1).Value = 'openai_gpt-oss-20b [T0-No-K]'
$ws.Cells.Item(34, 2).Value = '0.21 / 0.28'
$ws.Cells.Item(34, 3).Value = '0.38 / 0.46'
$ws.Cells.Item(34, 4).Value = '0.31 / 0.35'
$ws.Cells.Item(35, 1).Value = 'mistral-large-latest [T0.5-With-K]'
$ws.Cells.Item(35, 2).Value = '0.54 (0.49 - 0.57) / 0.67 (0.62 - 0.69)'
$ws.Cells.Item(35, 3).Value = '0.63 (0.61 - 0.67) / 0.71 (0.68 - 0.75)'
$ws.Cells.Item(35, 4).Value = '0.30 (0.29 - 0.32) / 0.42 (0.41 - 0.45)'
$ws.Cells.Item(36, 1).Value = 'gpt-4o-mini [T0.5-No-K]'
$ws.Cells.Item(36, 2).Value = '0.59 (0.55 - 0.62) / 0.65 (0.62 - 0.67)'
$ws.Cells.Item(36, 3).Value = '0.71 (0.69 - 0.73) / 0.82 (0.81 - 0.83)'
$ws.Cells.Item(36, 4).Value = '0.29 (0.28 - 0.31) / 0.35 (0.34 - 0.36)'
$ws.Cells.Item(37, 1).Value = 'deepseek-ai_DeepSeek-R1 [T0.5-With-K]'
$ws.Cells.Item(37, 2).Value = '0.39 (0.32 - 0.42) / 0.48 (0.42 - 0.52)'
$ws.Cells.Item(37, 3).Value = '0.45 (0.41 - 0.50) / 0.54 (0.49 - 0.60)'
$ws.Cells.Item(37, 4).Value = '0.28 (0.26 - 0.31) / 0.38 (0.37 - 0.38)'
$ws.Cells.Item(38, 1).Value = 'mistral-large-latest [T0.5-No-K]'
$ws.Cells.Item(38, 2).Value = '0.60 (0.59 - 0.61) / 0.71 (0.70 - 0.71)'
$ws.Cells.Item(38, 3).Value = '0.73 (0.72 - 0.75) / 0.82 (0.81 - 0.83)'
$ws.Cells.Item(38, 4).Value = '0.28 (0.25 - 0.34) / 0.40 (0.37 - 0.45)'
$ws.Cells.Item(39, 1).Value = 'gpt-4o-mini [T0.5-With-K]'
$ws.Cells.Item(39, 2).Value = '0.52 (0.49 - 0.56) / 0.57 (0.52 - 0.61)'
$ws.Cells.Item(39, 3).Value = '0.50 (0.43 - 0.54) / 0.61 (0.54 - 0.66)'
$ws.Cells.Item(39, 4).Value = '0.24 (0.23 - 0.25) / 0.31 (0.30 - 0.32)'
$ws.Cells.Item(40, 1).Value = 'gpt-4o-mini [T0-With-K]'
$ws.Cells.Item(40, 2).Value = '0.46 / 0.52'
$ws.Cells.Item(40, 3).Value = '0.47 / 0.60'
$ws.Cells.Item(40, 4).Value = '0.21 / 0.30'
$ws.Cells.Item(41, 1).Value = 'openai_gpt-oss-20b [T0-With-K]'
$ws.Cells.Item(41, 2).Value = '0.38 / 0.43'
$ws.Cells.Item(41, 3).Value = '0.20 / 0.27'
